$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (46075 -> 46076) for every data row from row 2 through row 236.
for ($r = 2; $r -le 236; $r++) {
    $ws.Cells.Item($r, 3).Value = 46076
}
